# Updated cryptos list - refresh Price / Volume(1h) figures, and update
# the ranking order of NEARProtocol / PaxDollar (rows 49 and 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.906.88"
$ws.Range("E2").Value = "  -3.70%  "

# Row 3
$ws.Range("D3").Value = "1.862.07"
$ws.Range("E3").Value = "  -2.90%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.42%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.58"
$ws.Range("E5").Value = "  -2.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4350"
$ws.Range("E7").Value = "  -5.39%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3725"
$ws.Range("E8").Value = "  -2.52%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07459"
$ws.Range("E9").Value = "  -3.44%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9327"
$ws.Range("E10").Value = "  -4.89%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.29"
$ws.Range("E11").Value = "  -4.29%  "

# Row 12
$ws.Range("D12").Value = "1.946.53"
$ws.Range("E12").Value = "  +2.97%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.713"
$ws.Range("E13").Value = "  -3.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.427"
$ws.Range("E14").Value = "  -4.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06867"
$ws.Range("E15").Value = "  -2.35%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.24%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.25"
$ws.Range("E17").Value = "  -3.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009037"
$ws.Range("E18").Value = "  -4.80%  "

# Row 19
$ws.Range("E19").Value = "  +0.16%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.84"
$ws.Range("E20").Value = "  -5.09%  "

# Row 21
$ws.Range("D21").Value = "27.903.92"
$ws.Range("E21").Value = "  -3.68%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.115"
$ws.Range("E22").Value = "  -4.13%  "

# Row 23
$ws.Range("E23").Value = "  +0.85%  "

# Row 24
$ws.Range("D24").Value = "2.155.07"
$ws.Range("E24").Value = "  +2.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.013"
$ws.Range("E25").Value = "  -3.43%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.72"
$ws.Range("E26").Value = "  -2.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.43"
$ws.Range("E27").Value = "  -3.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.491"
$ws.Range("E28").Value = "  -3.45%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.30"
$ws.Range("E29").Value = "  -4.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.711"
$ws.Range("E30").Value = "  -7.72%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09002"
$ws.Range("E31").Value = "  -3.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8146"
$ws.Range("E32").Value = "  -6.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.803"
$ws.Range("E33").Value = "  -5.99%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.173"
$ws.Range("E34").Value = "  -6.30%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.977"
$ws.Range("E35").Value = "  -2.51%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.003"
$ws.Range("E36").Value = "  +0.21%  "

# Row 37
$ws.Range("E37").Value = "  -2.80%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05494"
$ws.Range("E38").Value = "  -3.82%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01975"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.007"
$ws.Range("E40").Value = "  -1.29%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5254"
$ws.Range("E41").Value = "  -4.65%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.014"
$ws.Range("E42").Value = "  -6.60%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1699"
$ws.Range("E43").Value = "  -2.91%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.739"
$ws.Range("E44").Value = "  -6.83%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06766"
$ws.Range("E45").Value = "  -1.90%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4885"
$ws.Range("E46").Value = "  -5.70%  "

# Row 47
$ws.Range("E47").Value = "  -5.41%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.15"
$ws.Range("E48").Value = "  -2.99%  "

# --- Rows 49 & 50: NEARProtocol and PaxDollar swap ranking order ---
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.002"
$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.673"
$ws.Range("E50").Value = "  -6.03%  "

# --- Row 51 (only Volume changes) ---
$ws.Range("E51").Value = "  -15.64%  "
